# Weekly fruit/vegetable price update.
# A new weekly record (row) is inserted at row 29, pushing the existing
# rows 29-35 down to 30-36 (dimension grows from A1:R35 to A1:R36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 29, shifting rows 29-35 -> 30-36
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 44785
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100112035
$ws.Cells.Item(29, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 400
$ws.Cells.Item(29, 11).Value = 17000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 17425
$ws.Cells.Item(29, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 16).Value = 1162
$ws.Cells.Item(29, 17).Value = 15
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Make sure the newly inserted date cell keeps the date number format
# used by the rest of column D (style index 2 -> custom date format).
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(30, 4).NumberFormat
